# Mapping snippets to new API Documenter UID format
# Inserts a new "Member ID (methods only)" column into the Snippets table,
# between "Method/Prop/Rel Name" and "SnippetIdIntheYAMLFile", and marks
# the rows that describe methods (as opposed to properties/relationships)
# with a 1 in that new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# 1) Insert a blank worksheet column at C. This shifts the existing
#    SnippetIdIntheYAMLFile / MethodNameInTheSnippet columns one place to
#    the right (C->D, D->E) and carries over cell formatting from the
#    left-hand neighbour column into the freshly inserted cells.
$ws.Columns("C").Insert()

# 2) Grow the table definition so it covers the new column too.
$tbl.Resize($ws.Range("A1:E35"))

# 3) Re-assert the header text for every header cell in the table so the
#    table's column names line up with what's actually in row 1 (the
#    insert step above only moved values, it did not rename table
#    columns, so the old column names would otherwise end up "stuck" one
#    slot to the left of where they now live).
$ws.Range("A1").Value = "Class"
$ws.Range("B1").Value = "Method/Prop/Rel Name"
$ws.Range("C1").Value = "Member ID (methods only)"
$ws.Range("D1").Value = "SnippetIdIntheYAMLFile"
$ws.Range("E1").Value = "MethodNameInTheSnippet"

# 4) Give the new column the same width the author set for it (23, fixed
#    -- not an autofit "best fit" width like its neighbours).
$ws.Columns("C").ColumnWidth = 22.17

# 5) Fill in "1" for every row whose "Method/Prop/Rel Name" is an actual
#    method (insertBreak, search, set, ...) rather than a property or
#    relationship id. Rows not listed here keep the blank cell that the
#    column insert already produced for them.
$methodRows = @(5, 6, 7, 8, 9, 10, 11, 12, 14, 19, 25, 26, 27, 28, 29, 32, 33, 34)
foreach ($r in $methodRows) {
    $ws.Cells.Item($r, 3).Value = 1
}

# 6) Row 35 ("run" under Body.paragraphs) never had a cell in the old
#    layout's column C (no explicit style/content there), so after the
#    column insert it correctly has no cell in the new column C either --
#    nothing further to do for it.

# 7) Leave the final selection on C35, matching where the edit left off.
$ws.Range("C35").Select() | Out-Null
